$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '82.058.80'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '3.191.01'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'216.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.68%  '
$ws.Range("D6").Value = "'623.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.43%  '
$ws.Range("D7").Value = "'0.292"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +21.81%  '
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.586"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = '3.189.64'
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").Value = "'0.591"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.54%  '
$ws.Range("D12").Value = "'0.0000260"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +12.26%  '
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D14").Value = "'5.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.91%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.783.43'
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").Value = "'32.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '81.875.04'
$ws.Range("E17").Value = '  +2.87%  '
$ws.Range("D18").Value = '3.192.71'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = "'3.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.99%  '
$ws.Range("D20").Value = "'14.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.91%  '
$ws.Range("D21").Value = "'436.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("D22").Value = "'9.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.30%  '
$ws.Range("D23").Value = "'5.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = "'7.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.71%  '
$ws.Range("D25").Value = "'5.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +13.25%  '
$ws.Range("D26").Value = '3.363.37'
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").Value = "'11.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("D28").Value = "'76.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("E30").Value = '  +3.45%  '
$ws.Range("D31").Value = "'592.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.18%  '
$ws.Range("D32").Value = "'9.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("D35").Value = "'0.145"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +21.04%  '
$ws.Range("E36").Value = '  +9.20%  '
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").Value = "'22.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").Value = "'6.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.90%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = "'2.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.19%  '
$ws.Range("D43").Value = "'3.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +20.59%  '
$ws.Range("D44").Value = "'20.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.78%  '
$ws.Range("D45").Value = "'160.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.71%  '
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").Value = "'188.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.32%  '
$ws.Range("B48").Value = 'ImmutableX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D48").Value = "'1.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = "'44.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.36%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = "'26.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.71%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = "'0.773"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.95%  '
